$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 9.924897043638936
$ws.Range("B3").Value = 9.924897043638929
$ws.Range("B4").Value = 5.031000012360735
$ws.Range("B5").Value = 4.893897031278186
$ws.Range("B6").Value = 4.893897031278179
$ws.Range("B7").Value = 133.2838420986619
$ws.Range("B8").Value = 128.3899450673837
$ws.Range("B9").Value = 5.031000012360739
